# Cambios en diferentes materias
# Update the student's name on the title page.

$d = $word.ActiveDocument

$oldName = "Alma Nayeli Rodríguez Vázquez"
$newName = "Clara Margarita Fernández Riveron"

$found = $d.Content.Find.Execute($oldName, $true, $false, $false, $false, $false,
                                  $true, 1, $false, $newName, 2)

Write-Output "Replaced: $found"
